# Updating data file names.
#
# The underlying data edit swaps the "4:1 Unconditioned" (col D) and
# "1:4 Conditioned" (col E) columns -- header included -- for every row,
# plus a couple of cosmetic view/format tweaks left behind by Excel when
# the sheet was last saved (column E widened, view scrolled back to the
# top and zoom reset, selection moved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# Swap columns D and E (header row included) for every used row.
$dRange = $ws.Range("D1:D$lastRow")
$eRange = $ws.Range("E1:E$lastRow")

$dValues = $dRange.Value2
$eValues = $eRange.Value2

$dRange.Value = $eValues
$eRange.Value = $dValues

# Column E got a bit wider in the saved file.
$ws.Columns.Item(5).ColumnWidth = 16.75

# View cosmetics: zoom back to 100%, scroll to top, move the selection.
$win = $excel.ActiveWindow
$win.Zoom = 100
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("H5").Select()
